$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds version numbers stored as Text ("@") so they don't get
# reinterpreted as numbers/dates - apply that format (down through row 30)
# BEFORE writing any values, so new entries like "1.2.44" stay literal text
# instead of being parsed as a date.
$ws.Range("A1:A30").NumberFormat = "@"

# New release entry: version + description (goes into row 3)
$ws.Range("A3").Value = "1.2.44"
$ws.Range("B3").Value = "Добавлена связь .well-known/core с коммандным модулем"

# Columns B & C wrap their (often long) descriptions.
$ws.Range("B1:C4").WrapText = $true

# Header row and the new release row are taller to show the wrapped text.
$ws.Rows.Item(1).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# Printer/page setup for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
